$d = $word.ActiveDocument

# --- Locate the paragraph that needs the text rewrite ---
$old = "ElementOmhoog() and ElementOmlaag() methods are called by using the hoog or laag tools. They check the TekenElementLijst in reverse, to make sure they select the highest TekenElement visible. If there is a TekenElement selected, they will move it up or down the list respectively."
$new = "The ElementSelectie() method is called by using the hoog or laag tools. It checks the TekenElementLijst in reverse, to make sure it selects the newest TekenElement visible that has been hit. If there is a TekenElement selected, it will return it to the calling method."

$rng = $d.Content
$found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target sentence to replace."
}

$start = $rng.Start

# Replace the whole sentence (this also removes the separate trailing "." run)
$rng.Text = $new

$newEnd = $start + $new.Length

# --- Move the "_GoBack" bookmark to the end of the rewritten paragraph ---
# Re-adding a bookmark named "_GoBack" relocates Word's single reserved
# "last edit" bookmark, removing it from its previous location automatically.
#
# Placing a bookmark exactly at a paragraph-end position lands it incorrectly
# at the very start of the document, so a placeholder character is inserted
# right at the paragraph end first (moving the paragraph mark out of the way),
# the bookmark is anchored next to it, and the placeholder is removed again -
# the bookmark stays correctly anchored at the paragraph end.
$placeholder = $d.Range($newEnd, $newEnd)
$placeholder.InsertBefore("~")

$bmRange = $d.Range($newEnd, $newEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($newEnd, $newEnd + 1).Text = ""
